{"js": "// The source edit only touches build/tooling metadata that Word's\n// document model does not expose to automation:\n//   - the free-floating XML comment right after <w:body> in\n//     word/document.xml (\"... using REFERENCE JAXB in Microsoft Java\n//     21.0.8 ...\" -> \"... using REFERENCE JAXB in Oracle Java 21.0.8 ...\")\n//     is a docx4j/JAXB tool-stamp comment, not a content node, field,\n//     custom XML part, or document property, so it has no Word.* (or\n//     WordOpenXML) surface to read or write.\n//   - the reordering of the xmlns:* attributes on the <w:document>,\n//     <w:ftr>, <w:hdr> and <w:styles> root elements is just attribute\n//     serialization order from whatever tool re-saved the package; it\n//     carries no semantic meaning and Office.js has no API that lets a\n//     script control XML attribute order at all.\n//\n// In other words, every visible/editable piece of document content\n// (paragraphs, runs, formatting, headers, footers, styles, sections)\n// is identical before and after this change, so the correct,\n// content-faithful edit here is to leave the document body untouched.\n// We still touch the context so the script is a normal, valid Office.js\n// batch (load + sync) without mutating anything.\nconst body = context.document.body;\nbody.load(\"text\");\nawait context.sync();\n", "ps1": "# The source edit only touches build/tooling metadata that Word's COM\n# object model does not expose to automation:\n#   - the free-floating XML comment right after <w:body> in\n#     word/document.xml (\"... using REFERENCE JAXB in Microsoft Java\n#     21.0.8 ...\" -> \"... using REFERENCE JAXB in Oracle Java 21.0.8 ...\")\n#     is a docx4j/JAXB tool-stamp comment, not document text, a field, a\n#     custom XML part, or a document property, so there is no\n#     Range/Selection/Find/Content property that can see or change it\n#     (Content.Find.Execute(\"Microsoft Java\") matches nothing, because\n#     Word's Range text never includes XML comment nodes).\n#   - the reordering of the xmlns:* attributes on the <w:document>,\n#     <w:ftr>, <w:hdr> and <w:styles> root elements is just attribute\n#     serialization order from whichever tool re-saved the package; it\n#     has no semantic effect and there is no COM property that controls\n#     raw XML attribute order.\n#\n# Every visible/editable piece of document content (paragraphs, runs,\n# formatting, headers, footers, styles, sections) is identical before\n# and after this change, so the correct, content-faithful edit here is\n# to leave the document body untouched. We still touch the document so\n# this is a normal, valid COM script without mutating anything.\n$d = $word.ActiveDocument\n$null = $d.Content.Text\n"}
